$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Fix variable names: rows 21/22 (skinHorvath) and 25/26 (PedBE) were
# wrongly duplicating the "base" variable name instead of using the
# ageAcc2.* / ageAcc3.* naming convention.
$ws.Range("A25").Value = "ageAcc2.PedBE"
$ws.Range("A26").Value = "ageAcc3.PedBE"
$ws.Range("A21").Value = "ageAcc2.skinHorvath"
$ws.Range("A22").Value = "ageAcc3.skinHorvath"

$ws.Activate() | Out-Null
$ws.Range("A23").Select() | Out-Null
